$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$conv = $wb.Worksheets.Item("CONVERTION")
$lo = $ws.ListObjects.Item("Table1")

# --- Insert first new blank row at row 87 (pushes old 87.. down by one) ---
$ws.Rows(87).Insert()
$ws.Range("A88:K88").Copy()
$ws.Range("A87:K87").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A87").ClearContents()
$ws.Range("B87").ClearContents()
$ws.Range("C87").ClearContents()
$ws.Range("D87").ClearContents()
$ws.Range("E87").ClearContents()
$ws.Range("F87").ClearContents()
$ws.Range("G87").ClearContents()
$ws.Range("H87").ClearContents()
$ws.Range("I87").ClearContents()
$ws.Range("J87").ClearContents()
$ws.Range("K87").ClearContents()

# --- Insert second new blank row at row 90 (pushes old 89.. down by one more) ---
$ws.Rows(90).Insert()
$ws.Range("A89:K89").Copy()
$ws.Range("A90:K90").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A90").ClearContents()
$ws.Range("B90").ClearContents()
$ws.Range("C90").ClearContents()
$ws.Range("D90").ClearContents()
$ws.Range("E90").ClearContents()
$ws.Range("F90").ClearContents()
$ws.Range("G90").ClearContents()
$ws.Range("H90").ClearContents()
$ws.Range("I90").ClearContents()
$ws.Range("J90").ClearContents()
$ws.Range("K90").ClearContents()
$ws.Rows(90).RowHeight = 13.5

# --- Populate row 90 first so its PARTICULARS string ("UT(0-4-2)") becomes ---
# --- shared-string index 83, then row 87 ("UT(0-4-0)") becomes index 84,  ---
# --- matching the original author's edit order.                          ---
$ws.Range("B90").Value = "UT(0-4-2)"
$ws.Range("D90").Value = 0.504
$ws.Range("G90").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

$ws.Range("B87").Value = "UT(0-4-0)"
$ws.Range("D87").Value = 0.5
$ws.Range("G87").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- The row inserts above shifted the former last two data rows down to  ---
# --- rows 152/153; their calculated-column formula lost its structured   ---
# --- reference during the shift, so restore it explicitly.               ---
$ws.Range("G152").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("G153").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- Resize table to include the two new rows ---
$lo.Resize($ws.Range("A8:K153"))

# --- Update CONVERTION sheet minute/hour calculator inputs ---
$conv.Range("E3").Value = 4
$conv.Range("F3").ClearContents()

# --- Update selection in Sheet1 (cosmetic, matches author's last position) ---
$ws.Activate()
$ws.Range("F86").Select()

$wb.Application.CalculateFull()
